# Project deliverables completed !
#
# - Mark the remaining Phase 3 deliverables (rows 11, 14-17 on the
#   Checklist sheet) as fully completed and ready for printing.
# - Bump row 2's height to match the other (taller) data rows.
# - Leave the workbook scrolled back to the top of the Checklist sheet
#   with G1 selected, and land on the Legend tab as the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")
$legend = $wb.Worksheets.Item("Legend")

# Row 2 grew a bit taller (was using the sheet default height).
$ws.Rows.Item(2).RowHeight = 24

# Rows whose "Completeness" status moved from "Finalizing: 66-99 (%)"
# to fully "Completed: 100 (%)", and which are now flagged "Ready for
# printing" in column G.
$doneRows = @(11, 14, 15, 16, 17)

foreach ($r in $doneRows) {
    $statusCell = $ws.Cells.Item($r, 5)   # column E
    $readyCell = $ws.Cells.Item($r, 7)    # column G

    # A couple of these E cells (15 & 16) were still using the
    # borderless "Normal" style left over from before the table's
    # border formatting was applied to the rest of the column - align
    # them with a sibling cell's formatting first ...
    $ws.Cells.Item(13, 5).Copy()
    $statusCell.PasteSpecial(-4122)

    # ... then set the actual values.
    $statusCell.Value = "Completed: 100 (%)"
    $readyCell.Value = "Ready for printing"
}

$excel.CutCopyMode = 0

# Leave the Checklist view scrolled to the top with G1 selected ...
$ws.Range("G1").Select()

# ... and switch over to the Legend tab as the last thing before saving.
$legend.Activate()
